$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "KATEGORI"
$ws.Range("C1").Value = "SUBKATEGORI"

$ws.Cells.Item(1,1).EntireColumn.AutoFit() | Out-Null
$ws.Range("A1:N2").Columns.AutoFit() | Out-Null

$ws.Range("C2").Select()
